# "fehler der masken behoben" - fix the (16-bit) instruction masks in
# column D, which were missing a trailing hex digit (e.g. 0xff8 -> 0xff80,
# 0xff -> 0xff00, 0xfe -> 0xfe00, 0xfc -> 0xfc00, 0xf8 -> 0xf800).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D11").Value = "0xff80"
$ws.Range("D12").Value = "0xff80"
$ws.Range("D13").Value = "0xff80"

$ws.Range("D14").Value = "0xff00"
$ws.Range("D15").Value = "0xff00"
$ws.Range("D16").Value = "0xff00"
$ws.Range("D17").Value = "0xff00"
$ws.Range("D18").Value = "0xff00"
$ws.Range("D19").Value = "0xff00"
$ws.Range("D20").Value = "0xff00"
$ws.Range("D21").Value = "0xff00"
$ws.Range("D22").Value = "0xff00"
$ws.Range("D23").Value = "0xff00"
$ws.Range("D24").Value = "0xff00"
$ws.Range("D25").Value = "0xff00"
$ws.Range("D26").Value = "0xff00"
$ws.Range("D27").Value = "0xff00"
$ws.Range("D28").Value = "0xff00"
$ws.Range("D29").Value = "0xff00"
$ws.Range("D30").Value = "0xff00"

$ws.Range("D31").Value = "0xfe00"
$ws.Range("D32").Value = "0xfe00"

$ws.Range("D33").Value = "0xfc00"
$ws.Range("D34").Value = "0xfc00"
$ws.Range("D35").Value = "0xfc00"
$ws.Range("D36").Value = "0xfc00"
$ws.Range("D37").Value = "0xfc00"
$ws.Range("D38").Value = "0xfc00"

$ws.Range("D39").Value = "0xf800"
$ws.Range("D40").Value = "0xf800"

# View state: scrolled down to row 22, zoomed to 160%, single cell G27
# selected (was a multi-cell F27:J28 selection before).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("G27").Select()
$win.Zoom = 160
